# The presentation's design theme (ppt/theme/theme1.xml, "Integral") is
# swapped for the plain default "Office Theme" palette (the theme that
# previously only backed the notes master, ppt/theme/theme2.xml).
#
# PowerPoint exposes the twelve core theme colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink - in that fixed order) through
# Slide.ThemeColorScheme.Item(1..12).RGB, which writes straight into the
# slide master's theme part. Re-pointing each slot at the "Office Theme"
# values reproduces the colour swap the commit applied.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeColors = @(
    0,        # dk1      000000
    16777215, # lt1      FFFFFF
    6968388,  # dk2      44546A
    15132391, # lt2      E7E6E6
    13998939, # accent1  5B9BD5
    3243501,  # accent2  ED7D31
    10855845, # accent3  A5A5A5
    49407,    # accent4  FFC000
    12874308, # accent5  4472C4
    4697456,  # accent6  70AD47
    12673797, # hlink    0563C1
    7491477   # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
